# Auto-generated: scheduled market-data refresh for Leviathan Profits workbook.
# Updates currentAveragePrice* / Leve price / profit columns (H-N) for affected leve rows
# across all craft-job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 287.68182
$ws.Range("I11").Value = 287.68182
$ws.Range("K11").Value = 287.68182
$ws.Range("M11").Value = -147.68182
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").ClearContents()
$ws.Range("H64").Value = 4173.926
$ws.Range("I64").Value = 4009
$ws.Range("J64").Value = 4899.6
$ws.Range("K64").Value = 4009
$ws.Range("L64").Value = 4899.6
$ws.Range("M64").Value = -3761
$ws.Range("N64").Value = -5395.6
$ws.Range("H67").Value = 4173.926
$ws.Range("I67").Value = 4009
$ws.Range("J67").Value = 4899.6
$ws.Range("K67").Value = 4009
$ws.Range("L67").Value = 4899.6
$ws.Range("M67").Value = -3151
$ws.Range("N67").Value = -6615.6
$ws.Range("H113").Value = 4954.8887
$ws.Range("I113").Value = 4627.143
$ws.Range("K113").Value = 4627.143
$ws.Range("M113").Value = -1373.143
$ws.Range("H137").Value = 2282.2778
$ws.Range("I137").Value = 2003.8148
$ws.Range("J137").Value = 3117.6667
$ws.Range("K137").Value = 6011.4444
$ws.Range("L137").Value = 9353.000100000001
$ws.Range("M137").Value = -3461.4444
$ws.Range("N137").Value = -14453.0001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 41213.29
$ws.Range("I32").Value = 25219.887
$ws.Range("K32").Value = 25219.887
$ws.Range("M32").Value = -24932.887
$ws.Range("H61").Value = 1385.9584
$ws.Range("I61").Value = 1228.826
$ws.Range("K61").Value = 1228.826
$ws.Range("M61").Value = -1016.826
$ws.Range("H63").Value = 3000
$ws.Range("I63").Value = 3000
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 3000
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -2314
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 3000
$ws.Range("I66").Value = 3000
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 15000
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -11568
$ws.Range("N66").ClearContents()
$ws.Range("H74").Value = 1295.1552
$ws.Range("I74").Value = 1252.9615
$ws.Range("J74").Value = 1660.8334
$ws.Range("K74").Value = 1252.9615
$ws.Range("L74").Value = 1660.8334
$ws.Range("M74").Value = -378.9614999999999
$ws.Range("N74").Value = -3408.8334
$ws.Range("H77").Value = 1295.1552
$ws.Range("I77").Value = 1252.9615
$ws.Range("J77").Value = 1660.8334
$ws.Range("K77").Value = 6264.807499999999
$ws.Range("L77").Value = 8304.166999999999
$ws.Range("M77").Value = -1896.807499999999
$ws.Range("N77").Value = -17040.167
$ws.Range("H122").Value = 1167.2222
$ws.Range("I122").Value = 1167.2222
$ws.Range("K122").Value = 3501.6666
$ws.Range("M122").Value = -1051.6666
$ws.Range("H132").Value = 2258.743
$ws.Range("I132").Value = 1275.8518
$ws.Range("J132").Value = 5576
$ws.Range("K132").Value = 3827.5554
$ws.Range("L132").Value = 16728
$ws.Range("M132").Value = -1297.5554
$ws.Range("N132").Value = -21788
$ws.Range("H136").Value = 1385.9584
$ws.Range("I136").Value = 1228.826
$ws.Range("K136").Value = 3686.478
$ws.Range("M136").Value = -1136.478

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 7415.0513
$ws.Range("I20").Value = 6455.087
$ws.Range("J20").Value = 8795
$ws.Range("K20").Value = 6455.087
$ws.Range("L20").Value = 8795
$ws.Range("M20").Value = -6208.087
$ws.Range("N20").Value = -9289
$ws.Range("H35").Value = 35000
$ws.Range("J35").Value = 35000
$ws.Range("L35").Value = 35000
$ws.Range("N35").Value = -35620
$ws.Range("H82").Value = 18751
$ws.Range("J82").Value = 25000
$ws.Range("L82").Value = 25000
$ws.Range("N82").Value = -25766
$ws.Range("H85").Value = 18751
$ws.Range("J85").Value = 25000
$ws.Range("L85").Value = 25000
$ws.Range("N85").Value = -27652

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3739.7144
$ws.Range("I62").Value = 4000.7273
$ws.Range("J62").Value = 2782.6667
$ws.Range("K62").Value = 4000.7273
$ws.Range("L62").Value = 2782.6667
$ws.Range("M62").Value = -3376.7273
$ws.Range("N62").Value = -4030.6667
$ws.Range("H65").Value = 3739.7144
$ws.Range("I65").Value = 4000.7273
$ws.Range("J65").Value = 2782.6667
$ws.Range("K65").Value = 20003.6365
$ws.Range("L65").Value = 13913.3335
$ws.Range("M65").Value = -16883.6365
$ws.Range("N65").Value = -20153.3335
$ws.Range("H99").Value = 13867.8125
$ws.Range("I99").Value = 20006.625
$ws.Range("K99").Value = 20006.625
$ws.Range("M99").Value = -18508.625
$ws.Range("H122").Value = 2767.818
$ws.Range("I122").Value = 1815.4
$ws.Range("K122").Value = 5446.200000000001
$ws.Range("M122").Value = -2996.200000000001
$ws.Range("H126").Value = 13867.8125
$ws.Range("I126").Value = 20006.625
$ws.Range("K126").Value = 60019.875
$ws.Range("M126").Value = -57549.875
$ws.Range("H132").Value = 2041.7693
$ws.Range("I132").Value = 1962
$ws.Range("K132").Value = 5886
$ws.Range("M132").Value = -3356

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 674.75
$ws.Range("I113").Value = 479.8
$ws.Range("K113").Value = 1439.4
$ws.Range("M113").Value = 730.5999999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 12500
$ws.Range("J92").Value = 12500
$ws.Range("L92").Value = 12500
$ws.Range("N92").Value = -16244
$ws.Range("H122").Value = 3265.5
$ws.Range("I122").Value = 3444.5417
$ws.Range("J122").Value = 2549.3333
$ws.Range("K122").Value = 10333.6251
$ws.Range("L122").Value = 7647.999899999999
$ws.Range("M122").Value = -7883.625100000001
$ws.Range("N122").Value = -12547.9999
$ws.Range("H126").Value = 3850.0557
$ws.Range("I126").Value = 3706.4443
$ws.Range("J126").Value = 3993.6667
$ws.Range("K126").Value = 11119.3329
$ws.Range("L126").Value = 11981.0001
$ws.Range("M126").Value = -8649.332900000001
$ws.Range("N126").Value = -16921.0001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 46915
$ws.Range("J7").Value = 6000
$ws.Range("L7").Value = 6000
$ws.Range("N7").Value = -6224
$ws.Range("H22").Value = 970.25
$ws.Range("I22").Value = 698
$ws.Range("K22").Value = 698
$ws.Range("M22").Value = -403
$ws.Range("H27").Value = 970.25
$ws.Range("I27").Value = 698
$ws.Range("K27").Value = 698
$ws.Range("M27").Value = -591
$ws.Range("H122").Value = 3164.16
$ws.Range("I122").Value = 3026.261
$ws.Range("K122").Value = 9078.782999999999
$ws.Range("M122").Value = -6628.782999999999
$ws.Range("H126").Value = 46915
$ws.Range("J126").Value = 6000
$ws.Range("L126").Value = 18000
$ws.Range("N126").Value = -22940
$ws.Range("H136").Value = 2148.7646
$ws.Range("I136").Value = 1967.0714
$ws.Range("J136").Value = 2996.6667
$ws.Range("K136").Value = 5901.2142
$ws.Range("L136").Value = 8990.000100000001
$ws.Range("M136").Value = -3351.2142
$ws.Range("N136").Value = -14090.0001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 19984
$ws.Range("H73").Value = 19984
$ws.Range("H81").Value = 3856.5
$ws.Range("I81").Value = 3888.2222
$ws.Range("J81").Value = 3000
$ws.Range("K81").Value = 7776.4444
$ws.Range("L81").Value = 6000
$ws.Range("M81").Value = -6715.4444
$ws.Range("N81").Value = -8122
$ws.Range("H84").Value = 3856.5
$ws.Range("I84").Value = 3888.2222
$ws.Range("J84").Value = 3000
$ws.Range("K84").Value = 38882.222
$ws.Range("L84").Value = 30000
$ws.Range("M84").Value = -33578.222
$ws.Range("N84").Value = -40608
$ws.Range("H122").Value = 828.3333
$ws.Range("I122").Value = 828.3333
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 2484.9999
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -34.9998999999998
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 5623.552
$ws.Range("I132").Value = 6555.913
$ws.Range("J132").Value = 2049.5
$ws.Range("K132").Value = 19667.739
$ws.Range("L132").Value = 6148.5
$ws.Range("M132").Value = -17137.739
$ws.Range("N132").Value = -11208.5
